$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Sheet1 (Hoja1): add rows 4-7 ---
# Write column A first (Mariano, Paloma, Jose Manuel, Angel), then column E
# (IA, PR, SI, ISI) so new shared-string entries are created in the same
# order Excel originally produced them.
$ws1.Range("A4").Value = "Mariano"
$ws1.Range("A5").Value = "Paloma"
$ws1.Range("A6").Value = "Jose Manuel"
$ws1.Range("A7").Value = "Angel"

$ws1.Range("E4").Value = "IA"
$ws1.Range("E5").Value = "PR"
$ws1.Range("E6").Value = "SI"
$ws1.Range("E7").Value = "ISI"

$ws1.Range("F4").Value = 3
$ws1.Range("G4").Value = "Informatica"
$ws1.Range("H4").Value = "Mariano"

$ws1.Range("F5").Value = 3
$ws1.Range("G5").Value = "Informatica"
$ws1.Range("H5").Value = "Jose Manuel"

$ws1.Range("F6").Value = 3
$ws1.Range("G6").Value = "Informatica"
$ws1.Range("H6").Value = "Angel"

$ws1.Range("F7").Value = 3
$ws1.Range("G7").Value = "Informatica"
$ws1.Range("H7").Value = "Paloma"

# --- Sheet2 (Hoja2): update rows 4-5, add rows 6-14 ---
$ws2.Range("C4").Value = "Informatica 3"
$ws2.Range("C5").Value = "Informatica 3"

$ws2.Range("A6").Value = "IA"
$ws2.Range("B6").Value = 2
$ws2.Range("C6").Value = "Informatica 3"

$ws2.Range("A7").Value = "IA"
$ws2.Range("B7").Value = 2
$ws2.Range("C7").Value = "Informatica 3"

$ws2.Range("A8").Value = "IA"
$ws2.Range("B8").Value = 2
$ws2.Range("C8").Value = "Informatica 3"

$ws2.Range("A9").Value = "PRO"
$ws2.Range("B9").Value = 2
$ws2.Range("C9").Value = "Informatica 3"

$ws2.Range("A10").Value = "PRO"
$ws2.Range("B10").Value = 2
$ws2.Range("C10").Value = "Informatica 3"

$ws2.Range("A11").Value = "SI"
$ws2.Range("B11").Value = 2
$ws2.Range("C11").Value = "Informatica 3"

$ws2.Range("A12").Value = "SI"
$ws2.Range("B12").Value = 2
$ws2.Range("C12").Value = "Informatica 3"

$ws2.Range("A13").Value = "ISI"
$ws2.Range("B13").Value = 2
$ws2.Range("C13").Value = "Informatica 3"

$ws2.Range("A14").Value = "ISI"
$ws2.Range("B14").Value = 2
$ws2.Range("C14").Value = "Informatica 3"

$ws2.Range("A14").Select()

# Re-select Hoja1 last so it stays the active/tabSelected sheet, matching
# the original workbook state, while Hoja2 keeps its own remembered
# selection (A14) independently.
$ws1.Range("I7").Select()
